# Generate Report for handoff
#
# Updates the "Latest Handoff Datetime" (column D) for the row belonging to
# file "91b215e5-9d8f-4f58-97b9-d77ce90f9020.md" (row 5) on both the
# "zh-cn" and "de-de" localization-status worksheets, recording the new
# handoff timestamps produced while generating the handoff report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-16 03:52:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-16 03:53:02"
